# Daily attendance processing - reorder "Recorded By" (column G) entries so
# that a leading "System" entry is moved to the end of the comma-separated
# list instead of the front, e.g.:
#   "System, dnasr281@gmail.com"            -> "dnasr281@gmail.com, System"
#   "System, backup@backdoor.com, system"   -> "backup@backdoor.com, system, System"
# Rows whose "Recorded By" value does not start with "System, " (e.g. it
# already starts with another user, or is just "System" on its own) are left
# untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.UsedRange.Rows.Count
$prefix = "System, "
$updated = 0

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2
    if ($val -ne $null) {
        $text = $val.ToString()
        if ($text.StartsWith($prefix)) {
            $rest = $text.Substring($prefix.Length)
            $newVal = $rest + ", System"
            $cell.Value = $newVal
            $updated = $updated + 1
        }
    }
}

"Reordered Recorded By for $updated row(s)"
